$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 25
$ws.Range("A6").Value = 22
$ws.Range("A7").Value = 20
$ws.Range("A8").Value = 23
$ws.Range("A9").Value = 21
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = ' bwd packet length mean'
$ws.Range("A11").Value = 27
$ws.Range("B11").Value = ' bwd packets/s'
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = ' destination port'
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = ' flow duration'
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = ' flow iat max'
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = ' flow iat mean'
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = ' flow iat min'
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = ' flow iat std'
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = ' flow packets/s'
$ws.Range("A19").Value = 24
$ws.Range("B19").Value = ' fwd header length'
$ws.Range("A20").Value = 28
$ws.Range("B20").Value = ' fwd header length.1'
$ws.Range("A21").Value = 17
$ws.Range("B21").Value = ' fwd iat max'
$ws.Range("A22").Value = 15
$ws.Range("B22").Value = ' fwd iat mean'
$ws.Range("A23").Value = 18
$ws.Range("B23").Value = ' fwd iat min'
$ws.Range("A24").Value = 16
$ws.Range("B24").Value = ' fwd iat std'
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = ' fwd packet length max'
$ws.Range("A26").Value = 39
$ws.Range("B26").Value = ' idle max'
$ws.Range("A27").Value = 38
$ws.Range("B27").Value = ' idle std'
$ws.Range("B28").Value = ' init_win_bytes_backward'
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = ' protocol'
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = ' source port'
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = ' subflow bwd packets'
$ws.Range("A32").Value = 6
$ws.Range("B32").Value = ' total backward packets'
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = ' total fwd packets'
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 'active mean'
$ws.Range("A35").Value = 19
$ws.Range("B35").Value = 'bwd iat total'
$ws.Range("A36").Value = 14
$ws.Range("A37").Value = 26
$ws.Range("A39").Value = 31
$ws.Range("A40").Value = 29
